$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-12-22 Monday" "2025-12-23 Tuesday"
Replace-Text "332÷3=110, 2" "390÷5=78, 0"
Replace-Text "707÷8=88, 3" "328÷3=109, 1"
Replace-Text "998÷3=332, 2" "766÷3=255, 1"
Replace-Text "118÷8=14, 6" "266÷2=133, 0"
Replace-Text "640÷6=106, 4" "280÷7=40, 0"
Replace-Text "937÷9=104, 1" "143÷8=17, 7"
Replace-Text "366÷2=183, 0" "462÷8=57, 6"
Replace-Text "908÷4=227, 0" "320÷9=35, 5"
Replace-Text "109÷2=54, 1" "594÷3=198, 0"
Replace-Text "920÷5=184, 0" "362÷2=181, 0"
Replace-Text "132÷6=22, 0" "400÷4=100, 0"
Replace-Text "198÷5=39, 3" "850÷3=283, 1"
Replace-Text "646÷6=107, 4" "361÷2=180, 1"
Replace-Text "681÷7=97, 2" "288÷6=48, 0"
Replace-Text "601÷4=150, 1" "906÷8=113, 2"
Replace-Text "622÷9=69, 1" "183÷8=22, 7"
Replace-Text "706÷5=141, 1" "841÷6=140, 1"
Replace-Text "765÷3=255, 0" "826÷7=118, 0"
Replace-Text "254÷9=28, 2" "409÷9=45, 4"
Replace-Text "912÷6=152, 0" "364÷6=60, 4"
Replace-Text "103÷2=51, 1" "993÷3=331, 0"
Replace-Text "310÷4=77, 2" "951÷8=118, 7"
Replace-Text "308÷4=77, 0" "655÷9=72, 7"
Replace-Text "931÷5=186, 1" "619÷3=206, 1"
Replace-Text "214÷2=107, 0" "982÷5=196, 2"

Write-Host "Done"
